# filtered Km calculation for exponential growth phase only
# Updates the n_cbh_mean_se (column F) and n_cbh_median_se (column H) values
# for rows 2-20 to reflect the recalculated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.172905782177114
$ws.Range("H2").Value = 0.479568541489103

$ws.Range("F3").Value = 0.200262315313081
$ws.Range("H3").Value = 0.457945099017026

$ws.Range("F4").Value = 0.25022116058967
$ws.Range("H4").Value = 0.540270202668898

$ws.Range("F5").Value = 0.268341012420777
$ws.Range("H5").Value = 0.498401850322011

$ws.Range("F6").Value = 0.308198658675583
$ws.Range("H6").Value = 0.490220681240925

$ws.Range("F7").Value = 0.290216547300712
$ws.Range("H7").Value = 0.399285836436981

$ws.Range("F8").Value = 0.289015053142963
$ws.Range("H8").Value = 0.49241595574623

$ws.Range("F9").Value = 0.278280443974346
$ws.Range("H9").Value = 0.470235176397004

$ws.Range("F10").Value = 0.473212165263482
$ws.Range("H10").Value = 0.551546193172615

$ws.Range("F11").Value = 0.182520216912646
$ws.Range("H11").Value = 0.0499348925451281

$ws.Range("F12").Value = 0.186189592176234
$ws.Range("H12").Value = 0.475337454547499

$ws.Range("F13").Value = 0.293061786999847
$ws.Range("H13").Value = 0.45567932106213

$ws.Range("F14").Value = 0.417910060029995
$ws.Range("H14").Value = 0.443094044568794

$ws.Range("F15").Value = 0.44673037387372
$ws.Range("H15").Value = 0.721054033286632

$ws.Range("F16").Value = 1.107847845845
$ws.Range("H16").Value = 1.35901652893772

$ws.Range("F17").Value = 2.24732400026869
$ws.Range("H17").Value = 3.11507539149451

$ws.Range("F18").Value = 1.7015621600056
$ws.Range("H18").Value = 1.42348210431368

$ws.Range("F19").Value = 1.61298663182518
$ws.Range("H19").Value = 2.01718716138192

$ws.Range("F20").Value = 1.49476448346593
$ws.Range("H20").Value = 1.75850222635062
